$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.349.72"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "2.014.27"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.14"
$ws.Range("E5").Value = "  +5.27%  "

$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.99"
$ws.Range("E8").Value = "  -6.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  -3.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -4.73%  "

$ws.Range("E11").Value = "  -3.08%  "

$ws.Range("E12").Value = "  -6.75%  "

$ws.Range("D13").Value = "2.312.13"
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.60"
$ws.Range("E14").Value = "  -3.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.798"
$ws.Range("E15").Value = "  -7.63%  "

$ws.Range("E16").Value = "  -5.51%  "

$ws.Range("D17").Value = "2.008.00"
$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("D18").Value = "37.313.69"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.14"
$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("D20").Value = "0.0₃0838"
$ws.Range("E20").Value = "  -3.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "234.69"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.12"
$ws.Range("E22").Value = "  -2.91%  "

$ws.Range("E23").Value = "  +3.47%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.93"
$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.99"
$ws.Range("E27").Value = "  -5.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.62"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("E29").Value = "  -5.88%  "

$ws.Range("E30").Value = "  -4.97%  "

$ws.Range("E31").Value = "  -2.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.62"
$ws.Range("E32").Value = "  -4.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0641"
$ws.Range("E33").Value = "  -4.68%  "

$ws.Range("E34").Value = "  -0.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  -6.28%  "

$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.37"
$ws.Range("E38").Value = "  -6.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.51"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("E40").Value = "  +3.37%  "

$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("E43").Value = "  -5.59%  "

$ws.Range("D44").Value = "1.437.74"
$ws.Range("E44").Value = "  +3.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.81"
$ws.Range("E45").Value = "  -8.27%  "

$ws.Range("E46").Value = "  -3.64%  "

$ws.Range("E47").Value = "  -3.65%  "

$ws.Range("E48").Value = "  +2.42%  "

$ws.Range("E49").Value = "  -6.73%  "

$ws.Range("D50").Value = "2.203.63"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("E51").Value = "  -10.01%  "
